$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 46/47 swap: Aave <-> ARBITRUM (name, link, price, volume) ---
$ws.Range("B46").Value = "Aave"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "36.444.89"
$ws.Range("D3").Value = "1.948.21"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.10"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.46"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D12").Value = "2.235.94"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.826"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.01"
$ws.Range("D17").Value = "1.959.66"
$ws.Range("D18").Value = "36.388.82"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.45"
$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.93"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.05"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.81"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.135"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.22"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.39"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.23"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.21"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0968"
$ws.Range("D44").Value = "1.361.67"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.08"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.02"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.12"
$ws.Range("D50").Value = "2.126.05"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +12.97%  "
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -14.54%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.46%  "
